$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 415.5
$ws.Cells.Item(9, 9).Value = 311.33334
$ws.Cells.Item(9, 10).Value = 728
$ws.Cells.Item(9, 11).Value = 311.33334
$ws.Cells.Item(9, 12).Value = 728
$ws.Cells.Item(9, 13).Value = -142.33334
$ws.Cells.Item(9, 14).Value = -1066

$ws.Cells.Item(18, 8).Value = 6945144
$ws.Cells.Item(18, 9).Value = 6945144
$ws.Cells.Item(18, 11).Value = 6945144
$ws.Cells.Item(18, 13).Value = -6944860

$ws.Cells.Item(40, 8).Value = 1000
$ws.Cells.Item(40, 9).Value = 1000
$ws.Cells.Item(40, 11).Value = 1000
$ws.Cells.Item(40, 13).Value = -825

$ws.Cells.Item(43, 8).Value = 5396.5
$ws.Cells.Item(43, 10).Value = 5396.5
$ws.Cells.Item(43, 12).Value = 5396.5
$ws.Cells.Item(43, 14).Value = -5534.5

$ws.Cells.Item(134, 8).Value = 30000
$ws.Cells.Item(134, 10).Value = 30000
$ws.Cells.Item(134, 12).Value = 30000
$ws.Cells.Item(134, 14).Value = -40140

$ws.Cells.Item(137, 8).Value = 6458.923
$ws.Cells.Item(137, 9).Value = 6544.1816
$ws.Cells.Item(137, 10).Value = 5990
$ws.Cells.Item(137, 11).Value = 19632.5448
$ws.Cells.Item(137, 12).Value = 17970
$ws.Cells.Item(137, 13).Value = -17082.5448
$ws.Cells.Item(137, 14).Value = -23070

$ws.Cells.Item(141, 8).Value = 266236.88
$ws.Cells.Item(141, 9).Value = 1151.5807
$ws.Cells.Item(141, 10).Value = 853211.5
$ws.Cells.Item(141, 11).Value = 3454.7421
$ws.Cells.Item(141, 12).Value = 2559634.5
$ws.Cells.Item(141, 13).Value = 1725.2579
$ws.Cells.Item(141, 14).Value = -2569994.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 25002388
$ws.Cells.Item(2, 9).Value = 35715410
$ws.Cells.Item(2, 11).Value = 35715410
$ws.Cells.Item(2, 13).Value = -35715297

$ws.Cells.Item(61, 8).Value = 2762.1353
$ws.Cells.Item(61, 9).Value = 1048.1578
$ws.Cells.Item(61, 11).Value = 1048.1578
$ws.Cells.Item(61, 13).Value = -836.1578

$ws.Cells.Item(109, 8).Value = 25500
$ws.Cells.Item(109, 10).Value = 25500
$ws.Cells.Item(109, 12).Value = 25500
$ws.Cells.Item(109, 14).Value = -28274

$ws.Cells.Item(116, 8).Value = 25002388
$ws.Cells.Item(116, 9).Value = 35715410
$ws.Cells.Item(116, 11).Value = 35715410
$ws.Cells.Item(116, 13).Value = -35713116

$ws.Cells.Item(122, 8).Value = 3126.6667
$ws.Cells.Item(122, 9).Value = 2188.889
$ws.Cells.Item(122, 10).Value = 4533.3335
$ws.Cells.Item(122, 11).Value = 6566.667
$ws.Cells.Item(122, 12).Value = 13600.0005
$ws.Cells.Item(122, 13).Value = -4116.667
$ws.Cells.Item(122, 14).Value = -18500.0005

$ws.Cells.Item(136, 8).Value = 2762.1353
$ws.Cells.Item(136, 9).Value = 1048.1578
$ws.Cells.Item(136, 11).Value = 3144.4734
$ws.Cells.Item(136, 13).Value = -594.4733999999999

$ws.Cells.Item(137, 8).Value = 29642.857
$ws.Cells.Item(137, 10).Value = 29642.857
$ws.Cells.Item(137, 12).Value = 29642.857
$ws.Cells.Item(137, 14).Value = -39842.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 25002388
$ws.Cells.Item(3, 9).Value = 35715410
$ws.Cells.Item(3, 11).Value = 35715410
$ws.Cells.Item(3, 13).Value = -35715296

$ws.Cells.Item(22, 8).Value = 257.5
$ws.Cells.Item(22, 10).Value = 500
$ws.Cells.Item(22, 12).Value = 500
$ws.Cells.Item(22, 14).Value = -846

$ws.Cells.Item(94, 8).Value = 810.087
$ws.Cells.Item(94, 9).Value = 703.64703
$ws.Cells.Item(94, 11).Value = 703.64703
$ws.Cells.Item(94, 13).Value = -252.64703

$ws.Cells.Item(99, 8).Value = 4116.6665
$ws.Cells.Item(99, 9).Value = 563
$ws.Cells.Item(99, 10).Value = 7670.3335
$ws.Cells.Item(99, 11).Value = 563
$ws.Cells.Item(99, 12).Value = 7670.3335
$ws.Cells.Item(99, 13).Value = 935
$ws.Cells.Item(99, 14).Value = -10666.3335

$ws.Cells.Item(105, 8).Value = 2042.7222
$ws.Cells.Item(105, 9).Value = 2113.2222
$ws.Cells.Item(105, 10).Value = 1972.2222
$ws.Cells.Item(105, 11).Value = 2113.2222
$ws.Cells.Item(105, 12).Value = 1972.2222
$ws.Cells.Item(105, 13).Value = -366.2222000000002
$ws.Cells.Item(105, 14).Value = -5466.2222

$ws.Cells.Item(107, 8).Value = 3702.2
$ws.Cells.Item(107, 9).Value = 2061.8
$ws.Cells.Item(107, 10).Value = 5342.6
$ws.Cells.Item(107, 11).Value = 2061.8
$ws.Cells.Item(107, 12).Value = 5342.6
$ws.Cells.Item(107, 13).Value = -141.8000000000002
$ws.Cells.Item(107, 14).Value = -9182.6

$ws.Cells.Item(134, 8).Value = 2036.878
$ws.Cells.Item(134, 9).Value = 1615.0667
$ws.Cells.Item(134, 11).Value = 4845.2001
$ws.Cells.Item(134, 13).Value = -2310.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1328.7693
$ws.Cells.Item(22, 10).Value = 1639.9
$ws.Cells.Item(22, 12).Value = 1639.9
$ws.Cells.Item(22, 14).Value = -2339.9

$ws.Cells.Item(31, 8).Value = 8279.4
$ws.Cells.Item(31, 9).Value = 5500
$ws.Cells.Item(31, 10).Value = 8588.223
$ws.Cells.Item(31, 11).Value = 5500
$ws.Cells.Item(31, 12).Value = 8588.223
$ws.Cells.Item(31, 13).Value = -5205
$ws.Cells.Item(31, 14).Value = -9178.223

$ws.Cells.Item(34, 8).Value = 8279.4
$ws.Cells.Item(34, 9).Value = 5500
$ws.Cells.Item(34, 10).Value = 8588.223
$ws.Cells.Item(34, 11).Value = 5500
$ws.Cells.Item(34, 12).Value = 8588.223
$ws.Cells.Item(34, 13).Value = -5298
$ws.Cells.Item(34, 14).Value = -8992.223

$ws.Cells.Item(58, 8).Value = 9093033
$ws.Cells.Item(58, 9).Value = 1161.7715
$ws.Cells.Item(58, 10).Value = 25003806
$ws.Cells.Item(58, 11).Value = 1161.7715
$ws.Cells.Item(58, 12).Value = 25003806
$ws.Cells.Item(58, 13).Value = -958.7715000000001
$ws.Cells.Item(58, 14).Value = -25004212

$ws.Cells.Item(86, 8).Value = 6515.385
$ws.Cells.Item(86, 9).Value = 4950
$ws.Cells.Item(86, 10).Value = 7857.143
$ws.Cells.Item(86, 11).Value = 4950
$ws.Cells.Item(86, 12).Value = 7857.143
$ws.Cells.Item(86, 13).Value = -3827
$ws.Cells.Item(86, 14).Value = -10103.143

$ws.Cells.Item(89, 8).Value = 6515.385
$ws.Cells.Item(89, 9).Value = 4950
$ws.Cells.Item(89, 10).Value = 7857.143
$ws.Cells.Item(89, 11).Value = 24750
$ws.Cells.Item(89, 12).Value = 39285.715
$ws.Cells.Item(89, 13).Value = -19134
$ws.Cells.Item(89, 14).Value = -50517.715

$ws.Cells.Item(99, 8).Value = 2388.6667
$ws.Cells.Item(99, 9).Value = 1166.3334
$ws.Cells.Item(99, 11).Value = 1166.3334
$ws.Cells.Item(99, 13).Value = 331.6666

$ws.Cells.Item(126, 8).Value = 2388.6667
$ws.Cells.Item(126, 9).Value = 1166.3334
$ws.Cells.Item(126, 11).Value = 3499.0002
$ws.Cells.Item(126, 13).Value = -1029.0002

$ws.Cells.Item(132, 8).Value = 2201.196
$ws.Cells.Item(132, 9).Value = 1637.738
$ws.Cells.Item(132, 10).Value = 4830.6665
$ws.Cells.Item(132, 11).Value = 4913.214
$ws.Cells.Item(132, 12).Value = 14491.9995
$ws.Cells.Item(132, 13).Value = -2383.214
$ws.Cells.Item(132, 14).Value = -19551.9995

$ws.Cells.Item(134, 8).Value = 1839.3549
$ws.Cells.Item(134, 9).Value = 714.9524
$ws.Cells.Item(134, 10).Value = 4200.6
$ws.Cells.Item(134, 11).Value = 2144.8572
$ws.Cells.Item(134, 12).Value = 12601.8
$ws.Cells.Item(134, 13).Value = 390.1428000000001
$ws.Cells.Item(134, 14).Value = -17671.8

$ws.Cells.Item(136, 8).Value = 9093033
$ws.Cells.Item(136, 9).Value = 1161.7715
$ws.Cells.Item(136, 10).Value = 25003806
$ws.Cells.Item(136, 11).Value = 3485.3145
$ws.Cells.Item(136, 12).Value = 75011418
$ws.Cells.Item(136, 13).Value = -935.3145000000004
$ws.Cells.Item(136, 14).Value = -75016518

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 2337.5
$ws.Cells.Item(75, 9).Value = 1012.5
$ws.Cells.Item(75, 11).Value = 3037.5
$ws.Cells.Item(75, 13).Value = -2039.5

$ws.Cells.Item(78, 8).Value = 2337.5
$ws.Cells.Item(78, 9).Value = 1012.5
$ws.Cells.Item(78, 11).Value = 9112.5
$ws.Cells.Item(78, 13).Value = -4120.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 62018.293
$ws.Cells.Item(102, 9).Value = 2609.3333
$ws.Cells.Item(102, 10).Value = 204599.8
$ws.Cells.Item(102, 11).Value = 2609.3333
$ws.Cells.Item(102, 12).Value = 204599.8
$ws.Cells.Item(102, 13).Value = -987.3332999999998
$ws.Cells.Item(102, 14).Value = -207843.8

$ws.Cells.Item(113, 8).Value = 6000
$ws.Cells.Item(113, 9).Value = 4000
$ws.Cells.Item(113, 10).Value = 10000
$ws.Cells.Item(113, 11).Value = 4000
$ws.Cells.Item(113, 12).Value = 10000
$ws.Cells.Item(113, 13).Value = -1830
$ws.Cells.Item(113, 14).Value = -14340

$ws.Cells.Item(122, 8).Value = 4334
$ws.Cells.Item(122, 9).Value = 3239.2856
$ws.Cells.Item(122, 10).Value = 5727.273
$ws.Cells.Item(122, 11).Value = 9717.856800000001
$ws.Cells.Item(122, 12).Value = 17181.819
$ws.Cells.Item(122, 13).Value = -7267.856800000001
$ws.Cells.Item(122, 14).Value = -22081.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2163.6365
$ws.Cells.Item(40, 9).Value = 988.3333
$ws.Cells.Item(40, 10).Value = 3574
$ws.Cells.Item(40, 11).Value = 988.3333
$ws.Cells.Item(40, 12).Value = 3574
$ws.Cells.Item(40, 13).Value = -852.3333
$ws.Cells.Item(40, 14).Value = -3846

$ws.Cells.Item(132, 8).Value = 2837.5757
$ws.Cells.Item(132, 9).Value = 1414.2778
$ws.Cells.Item(132, 11).Value = 4242.8334
$ws.Cells.Item(132, 13).Value = -1712.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1727.7858
$ws.Cells.Item(113, 9).Value = 541.5714
$ws.Cells.Item(113, 10).Value = 2914
$ws.Cells.Item(113, 11).Value = 1624.7142
$ws.Cells.Item(113, 12).Value = 8742
$ws.Cells.Item(113, 13).Value = 545.2857999999999
$ws.Cells.Item(113, 14).Value = -13082

$ws.Cells.Item(122, 8).Value = 372636.94
$ws.Cells.Item(122, 9).Value = 501855.66
$ws.Cells.Item(122, 10).Value = 3440.5715
$ws.Cells.Item(122, 11).Value = 1505566.98
$ws.Cells.Item(122, 12).Value = 10321.7145
$ws.Cells.Item(122, 13).Value = -1503116.98
$ws.Cells.Item(122, 14).Value = -15221.7145

$ws.Cells.Item(132, 8).Value = 14797.768
$ws.Cells.Item(132, 9).Value = 3038.1177
$ws.Cells.Item(132, 11).Value = 9114.3531
$ws.Cells.Item(132, 13).Value = -6584.3531

$ws.Cells.Item(136, 8).Value = 927.1579
$ws.Cells.Item(136, 9).Value = 520.6875
$ws.Cells.Item(136, 11).Value = 1562.0625
$ws.Cells.Item(136, 13).Value = 987.9375
